$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 324 (Haba / Provincia de Limarí entry) and insert the copy
# above it. This shifts the existing row 324 (and everything below it,
# through row 363) down by one row, turning the sheet's data range from
# A1:R363 into A1:R364, while the new row 324 starts out as an exact copy
# of the original row 324.
$ws.Rows(324).Copy()
$ws.Rows(324).Insert()

# Now adjust the newly inserted row 324 so it represents the new weekly
# price observation: a later date (2023-07-20) and a higher sales volume,
# while keeping the rest of the row (quality, prices, unit, origin, etc.)
# the same as it was copied from the original row.
$ws.Range("D324").Value = "2023-07-20"
$ws.Range("J324").Value = 70
